# Auto-generated Excel COM-interop edit script
# Applies numeric 'want-to-go' (F) / price (G) count bumps plus the
# row 27-31 content reshuffle on sheet4 (全部类型) described by the diff.

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 919
$ws.Cells.Item(5, 6).Value = 1206
$ws.Cells.Item(6, 6).Value = 80
$ws.Cells.Item(7, 6).Value = 4494
$ws.Cells.Item(8, 6).Value = 2662
$ws.Cells.Item(10, 6).Value = 2619
$ws.Cells.Item(14, 6).Value = 1680
$ws.Cells.Item(15, 6).Value = 694
$ws.Cells.Item(16, 6).Value = 145
$ws.Cells.Item(17, 6).Value = 164
$ws.Cells.Item(18, 6).Value = 352
$ws.Cells.Item(20, 6).Value = 282
$ws.Cells.Item(21, 6).Value = 79
$ws.Cells.Item(22, 6).Value = 44
$ws.Cells.Item(24, 6).Value = 34
$ws.Cells.Item(25, 6).Value = 89
$ws.Cells.Item(26, 6).Value = 586
$ws.Cells.Item(27, 6).Value = 716
$ws.Cells.Item(28, 6).Value = 127
$ws.Cells.Item(30, 6).Value = 450
$ws.Cells.Item(31, 6).Value = 1636
$ws.Cells.Item(32, 6).Value = 1231
$ws.Cells.Item(33, 6).Value = 221
$ws.Cells.Item(34, 6).Value = 29
$ws.Cells.Item(35, 6).Value = 1292
$ws.Cells.Item(36, 6).Value = 2162
$ws.Cells.Item(37, 6).Value = 323
$ws.Cells.Item(39, 6).Value = 566
$ws.Cells.Item(40, 6).Value = 99
$ws.Cells.Item(41, 6).Value = 38
$ws.Cells.Item(43, 6).Value = 705
$ws.Cells.Item(44, 6).Value = 1385
$ws.Cells.Item(45, 6).Value = 144
$ws.Cells.Item(48, 6).Value = 63
$ws.Cells.Item(49, 6).Value = 86

# --- 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 14
$ws.Cells.Item(5, 6).Value = 77

# --- 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 919
$ws.Cells.Item(3, 6).Value = 14
$ws.Cells.Item(4, 6).Value = 80
$ws.Cells.Item(5, 6).Value = 4495
$ws.Cells.Item(6, 6).Value = 2662
$ws.Cells.Item(7, 6).Value = 2619
$ws.Cells.Item(8, 6).Value = 1680
$ws.Cells.Item(11, 6).Value = 694
$ws.Cells.Item(12, 6).Value = 145
$ws.Cells.Item(13, 6).Value = 164
$ws.Cells.Item(14, 6).Value = 352
$ws.Cells.Item(16, 6).Value = 282
$ws.Cells.Item(17, 6).Value = 79
$ws.Cells.Item(18, 6).Value = 44
$ws.Cells.Item(20, 6).Value = 34
$ws.Cells.Item(21, 6).Value = 586
$ws.Cells.Item(22, 6).Value = 716
$ws.Cells.Item(23, 6).Value = 127
$ws.Cells.Item(24, 6).Value = 77
$ws.Cells.Item(34, 6).Value = 2162
$ws.Cells.Item(35, 6).Value = 323
$ws.Cells.Item(39, 6).Value = 566
$ws.Cells.Item(40, 6).Value = 99
$ws.Cells.Item(41, 6).Value = 38
$ws.Cells.Item(43, 6).Value = 705
$ws.Cells.Item(44, 6).Value = 1385
$ws.Cells.Item(46, 6).Value = 144
$ws.Cells.Item(48, 6).Value = 86

# Rows 27-31: one event (2024-04-27 原神x崩铁周年特典only) dropped and
# a new ticket tier (配音演员紫枫儿内场票) inserted after row 30, shifting
# the intervening rows' content up by one.
# Row 27
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = '2024-05-01'
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).Value = '杭州·与梦回望动漫游戏展'
$ws.Cells.Item(27, 4).Value = '沈半路171号 T-Car杭州汽车文化主题公园'
$ws.Cells.Item(27, 5).Value = '2024.05.01 10:00-05.02 17:00'
$ws.Cells.Item(27, 6).Value = 450
$ws.Cells.Item(27, 7).Value = 70
$ws.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82725'
$ws.Cells.Item(27, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/lt13shal1710228931298.jpeg'

# Row 28
$ws.Cells.Item(28, 3).Value = '杭州·第37届 中二病 原神x星穹only'
$ws.Cells.Item(28, 4).Value = '康候圣街99号 顺丰创新中心'
$ws.Cells.Item(28, 5).Value = '2024.05.01 10:30-05.02 17:00'
$ws.Cells.Item(28, 6).Value = 1636
$ws.Cells.Item(28, 7).Value = 60
$ws.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82700'
$ws.Cells.Item(28, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/Kb75MESZ1710215541381.jpeg'

# Row 29
$ws.Cells.Item(29, 3).Value = '杭州·第7届YH樱花动漫游戏文化节'
$ws.Cells.Item(29, 4).Value = '德胜东路2539号 梦马汽车小镇'
$ws.Cells.Item(29, 5).Value = '2024.05.01 10:00-05.02 17:00'
$ws.Cells.Item(29, 6).Value = 1232
$ws.Cells.Item(29, 7).Value = 70
$ws.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82828'
$ws.Cells.Item(29, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/Kd0niodt1710905544733.jpeg'

# Row 30
$ws.Cells.Item(30, 3).Value = '杭州·第7届YH樱花漫展-SVIP嘉宾前排票'
$ws.Cells.Item(30, 6).Value = 221
$ws.Cells.Item(30, 7).Value = 168
$ws.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83267'
$ws.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/DgmIZ6G71711357279757.jpeg'

# Row 31
$ws.Cells.Item(31, 3).Value = '杭州·第7届YH樱花漫展-配音演员紫枫儿内场票'
$ws.Cells.Item(31, 5).Value = '2024.05.01 10:00-05.01 17:00'
$ws.Cells.Item(31, 6).Value = 29
$ws.Cells.Item(31, 7).Value = 98
$ws.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83331'
$ws.Cells.Item(31, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/h5ilz3SA1711351453471.jpeg'

